# feat: add 2022-Q1 data
#
# The old "总计" (totals) worksheet is repurposed to hold the new
# "2022-Q1" per-fund breakdown (it already carries the right header /
# column styling), and a brand-new "总计" worksheet is appended with the
# refreshed roll-up totals (2022-Q1 on top, 2020-Q4 below).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(2)
$totals.Name = "2022-Q1"

# Grab the existing styled header/index cells so the new columns can
# inherit the exact same formatting (bold font + border + centered).
$headerStyleSrc = $totals.Range("B1")
$indexStyleSrc = $totals.Range("A2")

# Update the 3 existing headers, add the 4 new ones.
$totals.Range("B1").Value = "基金代码"
$totals.Range("C1").Value = "基金名称"
$totals.Range("D1").Value = "基金规模"

$headerStyleSrc.Copy()
$totals.Range("E1:H1").PasteSpecial(-4122)
$totals.Range("E1").Value = "股票总仓位"
$totals.Range("F1").Value = "仓位占比"
$totals.Range("G1").Value = "持有市值(亿元)"
$totals.Range("H1").Value = "仓位排名"

# Data rows (2-6), column A carries the styled running index.
$indexStyleSrc.Copy()
$totals.Range("A3:A6").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4

# Columns B..G are plain text, so force text formatting before writing
# numeric-looking strings (fund codes, percentages, ...).
$totals.Range("B2:G6").NumberFormat = "@"

$totals.Range("B2").Value = "000988"
$totals.Range("C2").Value = "嘉实全球互联网股票 - 人民币QDII"
$totals.Range("D2").Value = "13.21"
$totals.Range("E2").Value = "85.88"
$totals.Range("F2").Value = "5.78"
$totals.Range("G2").Value = "0.7635"
$totals.Range("H2").Value = 7

$totals.Range("B3").Value = "000989"
$totals.Range("C3").Value = "嘉实全球互联网股票 - 美元现汇QDII"
$totals.Range("D3").Value = "13.21"
$totals.Range("E3").Value = "85.88"
$totals.Range("F3").Value = "5.78"
$totals.Range("G3").Value = "0.7635"
$totals.Range("H3").Value = 7

$totals.Range("B4").Value = "000990"
$totals.Range("C4").Value = "嘉实全球互联网股票 - 美元现钞QDII"
$totals.Range("D4").Value = "13.21"
$totals.Range("E4").Value = "85.88"
$totals.Range("F4").Value = "5.78"
$totals.Range("G4").Value = "0.7635"
$totals.Range("H4").Value = 7

$totals.Range("B5").Value = "006792"
$totals.Range("C5").Value = "鹏华香港美国互联网股票（LOF）美元现汇"
$totals.Range("D5").Value = "1.43"
$totals.Range("E5").Value = "83.72"
$totals.Range("F5").Value = "3.56"
$totals.Range("G5").Value = "0.0509"
$totals.Range("H5").Value = 6

$totals.Range("B6").Value = "160644"
$totals.Range("C6").Value = "鹏华香港美国互联网股票（LOF）人民币"
$totals.Range("D6").Value = "1.43"
$totals.Range("E6").Value = "83.72"
$totals.Range("F6").Value = "3.56"
$totals.Range("G6").Value = "0.0509"
$totals.Range("H6").Value = 6

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet with the updated roll-up totals
# ---------------------------------------------------------------------
$grand = $wb.Worksheets.Add($null, $totals)
$grand.Name = "总计"

$headerStyleSrc.Copy()
$grand.Range("B1:D1").PasteSpecial(-4122)
$grand.Range("B1").Value = "日期"
$grand.Range("C1").Value = "持有数量(只)"
$grand.Range("D1").Value = "持有市值(亿元)"

$indexStyleSrc.Copy()
$grand.Range("A2:A3").PasteSpecial(-4122)
$grand.Range("A2").Value = 0
$grand.Range("A3").Value = 1

$grand.Range("B2").Value = "2022-Q1"
$grand.Range("C2").Value = 5
$grand.Range("D2").Value = 2.39

$grand.Range("B3").Value = "2020-Q4"
$grand.Range("C3").Value = 1
$grand.Range("D3").Value = 0.02

# Keep "2020-Q4" as the active/selected tab, matching the untouched sheet1.
$wb.Worksheets.Item(1).Activate()
